$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update ticker in row 2: "AAAI" -> "24H" (the company name in B2 stays the same) ---
$ws.Range("A2").Value = "24H"

# --- Add two new rows (4 and 5) below the existing data, copying the look & feel of A3 ---
$ws.Range("A3").Copy($ws.Range("A4"))
$ws.Range("A3").Copy($ws.Range("A5"))
$ws.Range("A4").Value = "AACORP"
$ws.Range("A5").Value = "AAM"

# --- Rework hyperlinks: the old A2/A3 links are replaced, B2/B3 are left untouched ---
$links = @()
foreach ($h in $ws.Hyperlinks) { $links += $h }
$addrs = @()
foreach ($h in $links) { $addrs += $h.Range.Address() }
for ($i = $links.Count - 1; $i -ge 0; $i--) {
    if ($addrs[$i] -eq '$A$2' -or $addrs[$i] -eq '$A$3') {
        $links[$i].Delete()
    }
}

$ws.Hyperlinks.Add($ws.Range("A2"), "http://s.cafef.vn/otc/24H-ctcp-quang-cao-truc-tuyen-24h.chn", [Type]::Missing, [Type]::Missing, "http://s.cafef.vn/otc/24H-ctcp-quang-cao-truc-tuyen-24h.chn")
$ws.Hyperlinks.Add($ws.Range("A3"), "http://s.cafef.vn/hose/AAA-cong-ty-co-phan-nhua-va-moi-truong-xanh-an-phat.chn", [Type]::Missing, [Type]::Missing, "http://s.cafef.vn/hose/AAA-cong-ty-co-phan-nhua-va-moi-truong-xanh-an-phat.chn")
$ws.Hyperlinks.Add($ws.Range("A4"), "http://s.cafef.vn/otc/AACORP-ctcp-xay-dung-kien-truc-aa.chn", [Type]::Missing, [Type]::Missing, "http://s.cafef.vn/otc/AACORP-ctcp-xay-dung-kien-truc-aa.chn")
$ws.Hyperlinks.Add($ws.Range("A5"), "http://s.cafef.vn/hose/AAM-cong-ty-co-phan-thuy-san-mekong.chn", [Type]::Missing, [Type]::Missing, "http://s.cafef.vn/hose/AAM-cong-ty-co-phan-thuy-san-mekong.chn")

# --- Adding hyperlinks re-applies a slightly different style variant, and baking in
#     TextToDisplay overwrites the cell text; restore the original formatting (style 7
#     look) and the correct ticker text on every cell that now carries a hyperlink ---
$ws.Range("A2").Copy($ws.Range("A2"))
$ws.Range("A3").Copy($ws.Range("A3"))
$ws.Range("A3").Copy($ws.Range("A4"))
$ws.Range("A3").Copy($ws.Range("A5"))
$ws.Range("A2").Value = "24H"
$ws.Range("A3").Value = "AAA"
$ws.Range("A4").Value = "AACORP"
$ws.Range("A5").Value = "AAM"
